# The weekly/daily data rows (2..10) get re-shuffled: each row's Fecha (D),
# Volumen (M), Unidad de comercializacion (Q), Precio $/Kg (S) and Kg/unidad (T)
# move to a different row while the remaining descriptive columns stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 2..10

# Snapshot current values of the columns that get reshuffled.
$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @($ws.Cells.Item($r, 4).Value2, $ws.Cells.Item($r, 13).Value2, $ws.Cells.Item($r, 17).Value2, $ws.Cells.Item($r, 19).Value2, $ws.Cells.Item($r, 20).Value2)
}

# Mapping: target row -> source row that supplies the (D, M, Q, S, T) tuple.
$mapping = @{ 2 = 9; 3 = 2; 4 = 8; 5 = 4; 6 = 6; 7 = 5; 8 = 10; 9 = 7; 10 = 3 }

foreach ($r in $rows) {
    $src = $mapping[$r]
    $vals = $orig[$src]
    $ws.Cells.Item($r, 4).Value2 = $vals[0]
    $ws.Cells.Item($r, 13).Value2 = $vals[1]
    $ws.Cells.Item($r, 17).Value2 = $vals[2]
    $ws.Cells.Item($r, 19).Value2 = $vals[3]
    $ws.Cells.Item($r, 20).Value2 = $vals[4]
}
